# Common card now can have 2 copies in the deck.
# The "已经入册/同种卡片只能入册1张" error strings are replaced by
# "该卡牌只能入册1张" (1 copy) / "该卡牌只能入册2张" (2 copies).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Error")

# Error 1001 first, so its new string is appended to the shared-string
# table before the one for error 1000 (matches the author's save order).
$ws.Range("B12").Value = "该卡牌只能入册2张"
$ws.Range("B11").Value = "该卡牌只能入册1张"

# Leave the selection where the author left it when saving.
$ws.Range("B10").Select()
